$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-29 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-30 Monday", 2)

$d.Content.Find.Execute("450×6=2700", $true, $false, $false, $false, $false, $true, 1, $false, "595×9=5355", 2)
$d.Content.Find.Execute("647×7=4529", $true, $false, $false, $false, $false, $true, 1, $false, "152×6=912", 2)
$d.Content.Find.Execute("532×4=2128", $true, $false, $false, $false, $false, $true, 1, $false, "806×6=4836", 2)
$d.Content.Find.Execute("644×9=5796", $true, $false, $false, $false, $false, $true, 1, $false, "542×8=4336", 2)
$d.Content.Find.Execute("225×6=1350", $true, $false, $false, $false, $false, $true, 1, $false, "132×5=660", 2)

$d.Content.Find.Execute("768×4=3072", $true, $false, $false, $false, $false, $true, 1, $false, "294×4=1176", 2)
$d.Content.Find.Execute("454×4=1816", $true, $false, $false, $false, $false, $true, 1, $false, "857×4=3428", 2)
$d.Content.Find.Execute("282×6=1692", $true, $false, $false, $false, $false, $true, 1, $false, "161×5=805", 2)
$d.Content.Find.Execute("420×8=3360", $true, $false, $false, $false, $false, $true, 1, $false, "920×8=7360", 2)
$d.Content.Find.Execute("633×2=1266", $true, $false, $false, $false, $false, $true, 1, $false, "536×5=2680", 2)

$d.Content.Find.Execute("873×8=6984", $true, $false, $false, $false, $false, $true, 1, $false, "422×8=3376", 2)
$d.Content.Find.Execute("494×9=4446", $true, $false, $false, $false, $false, $true, 1, $false, "548×6=3288", 2)
$d.Content.Find.Execute("694×8=5552", $true, $false, $false, $false, $false, $true, 1, $false, "525×5=2625", 2)
$d.Content.Find.Execute("152×2=304", $true, $false, $false, $false, $false, $true, 1, $false, "635×5=3175", 2)
$d.Content.Find.Execute("889×2=1778", $true, $false, $false, $false, $false, $true, 1, $false, "493×5=2465", 2)

$d.Content.Find.Execute("721×5=3605", $true, $false, $false, $false, $false, $true, 1, $false, "469×2=938", 2)
$d.Content.Find.Execute("994×9=8946", $true, $false, $false, $false, $false, $true, 1, $false, "185×3=555", 2)
$d.Content.Find.Execute("338×5=1690", $true, $false, $false, $false, $false, $true, 1, $false, "187×8=1496", 2)
$d.Content.Find.Execute("681×6=4086", $true, $false, $false, $false, $false, $true, 1, $false, "432×9=3888", 2)
$d.Content.Find.Execute("810×8=6480", $true, $false, $false, $false, $false, $true, 1, $false, "400×7=2800", 2)

$d.Content.Find.Execute("499×2=998", $true, $false, $false, $false, $false, $true, 1, $false, "762×9=6858", 2)
$d.Content.Find.Execute("288×6=1728", $true, $false, $false, $false, $false, $true, 1, $false, "564×8=4512", 2)
$d.Content.Find.Execute("156×9=1404", $true, $false, $false, $false, $false, $true, 1, $false, "591×3=1773", 2)
$d.Content.Find.Execute("276×4=1104", $true, $false, $false, $false, $false, $true, 1, $false, "459×3=1377", 2)
$d.Content.Find.Execute("854×2=1708", $true, $false, $false, $false, $false, $true, 1, $false, "779×3=2337", 2)
